$d = $word.ActiveDocument

# Remove "Backbone.JS, " from the front-end skills line:
#   "...Javascript (jQuery, Backbone.JS, Angular.JS,)"
# becomes
#   "...Javascript (jQuery, Angular.JS,)"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Backbone.JS, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
